$wb = $excel.ActiveWorkbook

# --- Enrichment protocol sheet: row 4 "process_core" labels (previously
# duplicated under a separate "protocol_core" naming) are corrected to
# reuse the same strings as the header row (row 2). ---
$wsEnrichment = $wb.Worksheets.Item("Enrichment protocol")
$wsEnrichment.Range("B4").Value = "enrichment_protocol.process_core.protocol_name"
$wsEnrichment.Range("C4").Value = "enrichment_protocol.process_core.protocol_description"
$wsEnrichment.Range("D4").Value = "enrichment_protocol.process_core.start_time"
$wsEnrichment.Range("E4").Value = "enrichment_protocol.process_core.process_location"
$wsEnrichment.Range("F4").Value = "enrichment_protocol.process_core.operator_identity"
$wsEnrichment.Activate()
$wsEnrichment.Range("A6:XFD6").Select()

# --- Library preparation protocol sheet: same kind of fix. ---
$wsLibPrep = $wb.Worksheets.Item("Library preparation protocol")
$wsLibPrep.Range("A4").Value = "library_preparation_protocol.process_core.protocol_id"
$wsLibPrep.Range("B4").Value = "library_preparation_protocol.process_core.protocol_name"
$wsLibPrep.Range("C4").Value = "library_preparation_protocol.process_core.protocol_description"
$wsLibPrep.Range("D4").Value = "library_preparation_protocol.process_core.start_time"
$wsLibPrep.Range("E4").Value = "library_preparation_protocol.process_core.process_location"
$wsLibPrep.Range("F4").Value = "library_preparation_protocol.process_core.operator_identity"
$wsLibPrep.Activate()
$excel.ActiveWindow.ScrollColumn = 28
$wsLibPrep.Range("A6:XFD6").Select()

# --- Sequencing protocol sheet: same kind of fix, plus L4 realigned to
# the process_type string used elsewhere on the sheet. ---
$wsSeq = $wb.Worksheets.Item("Sequencing protocol")
$wsSeq.Range("A4").Value = "sequencing_protocol.process_core.process_id"
$wsSeq.Range("B4").Value = "sequencing_protocol.process_core.process_name"
$wsSeq.Range("C4").Value = "sequencing_protocol.process_core.process_description"
$wsSeq.Range("D4").Value = "sequencing_protocol.process_core.start_time"
$wsSeq.Range("E4").Value = "sequencing_protocol.process_core.process_location"
$wsSeq.Range("F4").Value = "sequencing_protocol.process_core.operator_identity"
$wsSeq.Range("L4").Value = "sequencing_protocol.process_type.text"
$wsSeq.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$wsSeq.Range("M14").Select()

# --- Sequence files sheet stays the active tab, with an updated selection. ---
$wsSeqFiles = $wb.Worksheets.Item("Sequence files")
$wsSeqFiles.Activate()
$wsSeqFiles.Range("L9").Select()
